# Derek's Log.xlsx - "Spell checked and permission changed"
# Adds a new FRIDAY section (rows 46-49) to the "Logs" sheet:
#   - row 46: section header "FRIDAY" (style copied from the existing
#     "MONDAY" header row, 43)
#   - row 47: Pickup Mic task in room 152 (trailing-space variant)
#   - row 48: AV Shutdown task in room 152
#   - row 49: Other (pickup wireless keyboard) task in room 152

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Row 46: new "FRIDAY" section header --------------------------------
# Copy formatting + structure from row 43 (the "MONDAY" header row) so the
# new header row picks up the same shaded/bordered style set, then
# overwrite its label.
$ws.Range("A43:F43").Copy($ws.Range("A46:F46"))
$ws.Range("C46").Value = "FRIDAY"

# --- Row 47: Pickup Mic ---------------------------------------------------
$ws.Range("A47").Value = "Pickup Mic"
$ws.Range("B47").Value = 42594
$ws.Range("C47").Value = "1600"
$ws.Range("D47").Value = "FC"
$ws.Range("E47").Value = "152 "
$ws.Range("F47").Value = "Pick up mic and clip and mic stand - return to Founders 156A storeroom. Leave mic cables and matts. Turn off amp. "

# --- Row 48: AV Shutdown ---------------------------------------------------
$ws.Range("A48").Value = "AV Shutdown"
$ws.Range("B48").Value = 42594
$ws.Range("C48").Value = "1600"
$ws.Range("D48").Value = "FC"
$ws.Range("E48").Value = "152"
$ws.Range("F48").Value = "Turn off PC and projector. Leave equipment in room. Lock room. Key for room in Founders 164 storeroom."

# --- Row 49: Other (wireless keyboard pickup) ------------------------------
$ws.Range("A49").Value = "Other"
$ws.Range("B49").Value = 42594
$ws.Range("C49").Value = "1600"
$ws.Range("D49").Value = "FC"
$ws.Range("E49").Value = "152"
$ws.Range("F49").Value = "Pick up wireless keyboard and return to Founders 156A storeroom. "

# The new rows wrap onto two lines in column F, same as the existing
# wrapped rows (34, 40, 42, 44) which render at height 30.
$ws.Rows.Item(47).RowHeight = 30
$ws.Rows.Item(48).RowHeight = 30
$ws.Rows.Item(49).RowHeight = 30

# Move the active selection the way the author's last save left it.
$ws.Range("F52").Select()
